$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 334 (a new weekly price observation),
# shifting all existing rows 334-355 down by one (to 335-356).
$ws.Rows.Item(334).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(334, 1).Value = 4
$ws.Cells.Item(334, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(334, 3).Value = "Los Lagos"
$ws.Cells.Item(334, 4).Value = 44931
$ws.Cells.Item(334, 5).Value = 10
$ws.Cells.Item(334, 6).Value = "Fruta"
$ws.Cells.Item(334, 7).Value = 100108
$ws.Cells.Item(334, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(334, 9).Value = 100108005
$ws.Cells.Item(334, 10).Value = "Piña"
$ws.Cells.Item(334, 11).Value = "Caramelo"
$ws.Cells.Item(334, 12).Value = "Segunda"
$ws.Cells.Item(334, 13).Value = 120
$ws.Cells.Item(334, 14).Value = 24000
$ws.Cells.Item(334, 15).Value = 24000
$ws.Cells.Item(334, 16).Value = 24000
$ws.Cells.Item(334, 17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(334, 18).Value = "Ecuador"
$ws.Cells.Item(334, 19).Value = 1714
$ws.Cells.Item(334, 20).Value = 14
